$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2079207920792079
$ws.Range("C2").Value = 0.5313531353135313
$ws.Range("J2").Value = 0.0165016501650165
$ws.Range("O2").Value = 0.009900990099009901
$ws.Range("P2").Value = 0.1353135313531353
$ws.Range("S2").Value = 0.09900990099009901

$ws.Range("B3").Value = 0.02424242424242424
$ws.Range("C3").Value = 0.02424242424242424
$ws.Range("J3").Value = 0.01818181818181818
$ws.Range("P3").Value = 0.7090909090909091
$ws.Range("S3").Value = 0.2242424242424242

$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("P4").Value = 0.5510204081632653
$ws.Range("S4").Value = 0.3469387755102041

$ws.Range("B6").Value = 0.05581395348837209
$ws.Range("D6").Value = 0.0186046511627907
$ws.Range("F6").Value = 0.05116279069767442
$ws.Range("J6").Value = 0.2279069767441861
$ws.Range("O6").Value = 0.05116279069767442
$ws.Range("Q6").Value = 0.2186046511627907
$ws.Range("R6").Value = 0.05116279069767442
$ws.Range("S6").Value = 0.3255813953488372

$ws.Range("B7").Value = 0.1256544502617801
$ws.Range("D7").Value = 0.03141361256544502
$ws.Range("F7").Value = 0.05235602094240838
$ws.Range("J7").Value = 0.193717277486911
$ws.Range("O7").Value = 0.005235602094240838
$ws.Range("Q7").Value = 0.1727748691099476
$ws.Range("R7").Value = 0.07329842931937172
$ws.Range("S7").Value = 0.3455497382198953

$ws.Range("B8").Value = 0.08874458874458875
$ws.Range("D8").Value = 0.02164502164502164
$ws.Range("E8").Value = 0.002164502164502165
$ws.Range("F8").Value = 0.06277056277056277
$ws.Range("O8").Value = 0.01948051948051948
$ws.Range("Q8").Value = 0.1601731601731602
$ws.Range("R8").Value = 0.06493506493506493
$ws.Range("S8").Value = 0.4632034632034632

$ws.Range("B9").Value = 0.09677419354838709
$ws.Range("D9").Value = 0.008064516129032258
$ws.Range("F9").Value = 0.07258064516129033
$ws.Range("J9").Value = 0.08870967741935484
$ws.Range("O9").Value = 0.01612903225806452
$ws.Range("Q9").Value = 0.1612903225806452
$ws.Range("R9").Value = 0.0564516129032258
$ws.Range("S9").Value = 0.5

$ws.Range("B10").Value = 0.1199671322925226
$ws.Range("D10").Value = 0.02465078060805259
$ws.Range("E10").Value = 0.0008216926869350862
$ws.Range("F10").Value = 0.07559572719802794
$ws.Range("J10").Value = 0.1027115858668858
$ws.Range("O10").Value = 0.02218570254724733
$ws.Range("Q10").Value = 0.1922760887428102
$ws.Range("R10").Value = 0.07641741988496302
$ws.Range("S10").Value = 0.3853738701725555

$ws.Range("F11").Value = 0.003144654088050315
$ws.Range("G11").Value = 0.119496855345912
$ws.Range("J11").Value = 0.1352201257861635
$ws.Range("K11").Value = 0.1886792452830189
$ws.Range("L11").Value = 0.5440251572327044
$ws.Range("S11").Value = 0.009433962264150943

$ws.Range("G12").Value = 0.7175141242937854
$ws.Range("J12").Value = 0.2090395480225989
$ws.Range("K12").Value = 0.02259887005649718
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("S12").Value = 0.01694915254237288

$ws.Range("G13").Value = 0.6444444444444445
$ws.Range("J13").Value = 0.3555555555555556

$ws.Range("F15").Value = 0.0136986301369863
$ws.Range("H15").Value = 0.1917808219178082
$ws.Range("I15").Value = 0.0547945205479452
$ws.Range("J15").Value = 0.3378995433789954
$ws.Range("K15").Value = 0.0593607305936073
$ws.Range("M15").Value = 0.0319634703196347
$ws.Range("O15").Value = 0.0684931506849315
$ws.Range("S15").Value = 0.2420091324200913

$ws.Range("F16").Value = 0.02222222222222222
$ws.Range("H16").Value = 0.1277777777777778
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.45
$ws.Range("K16").Value = 0.1222222222222222
$ws.Range("M16").Value = 0.01666666666666667
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.1555555555555556

$ws.Range("F17").Value = 0.02689486552567237
$ws.Range("H17").Value = 0.2200488997555012
$ws.Range("I17").Value = 0.05378973105134474
$ws.Range("J17").Value = 0.4058679706601467
$ws.Range("K17").Value = 0.1026894865525672
$ws.Range("M17").Value = 0.01711491442542787
$ws.Range("O17").Value = 0.05623471882640587
$ws.Range("S17").Value = 0.117359413202934

$ws.Range("F18").Value = 0.03821656050955414
$ws.Range("H18").Value = 0.1719745222929936
$ws.Range("I18").Value = 0.07006369426751592
$ws.Range("J18").Value = 0.4331210191082803
$ws.Range("K18").Value = 0.1019108280254777
$ws.Range("M18").Value = 0.01273885350318471
$ws.Range("O18").Value = 0.05732484076433121
$ws.Range("S18").Value = 0.1146496815286624

$ws.Range("F19").Value = 0.01271860095389507
$ws.Range("H19").Value = 0.2241653418124006
$ws.Range("I19").Value = 0.05564387917329094
$ws.Range("J19").Value = 0.3688394276629571
$ws.Range("K19").Value = 0.1271860095389507
$ws.Range("M19").Value = 0.02066772655007949
$ws.Range("N19").Value = 0.001589825119236884
$ws.Range("O19").Value = 0.06756756756756757
$ws.Range("S19").Value = 0.1216216216216216

